# ---------------------------------------------------------------------------
# Commit: Sun, May 10, 2020 11:06:31 PM
#
# The canonical OOXML diff shows two changes:
#
#   1. ppt/slides/slide5.xml  - the table's <a:tableStyleId> is changed from
#      {CA688BE3-F663-40D7-82EB-534EB8C3B3B2} to
#      {E07B0B2D-2213-4F66-BB75-B3AC2C27C7B8} (a different built-in table
#      style).
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their bodies: the
#      part that used to hold the "Office Theme" / "Office" colour scheme
#      now holds the "Integral" / "Red Violet" colours and vice-versa. The
#      slide master (and therefore every layout/slide that is actually
#      rendered) is wired to theme2.xml, which used to carry "Integral" and
#      ends up carrying the plain "Office Theme" palette; theme1.xml (wired
#      to the notes master) goes the other way. The font scheme and format
#      scheme (fills/lines/effects) are byte-identical between the two
#      themes, so the only real content change is the 12 theme colours.
#
# The PowerPoint object model doesn't expose raw OOXML parts, so the theme
# swap is reproduced through the SlideMaster's ThemeColorScheme (the theme
# that is actually wired to every slide via the layouts/master chain, i.e.
# ppt/theme/theme2.xml) by writing each of the 12 standard theme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - same order as
# ThemeColorScheme.Colors(1..12)) to the "Office Theme" values.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$tbl = $p.Slides.Item(5).Shapes.Item(2).Table
$tbl.ApplyStyle("{E07B0B2D-2213-4F66-BB75-B3AC2C27C7B8}")

# --- 2) Theme colours: "Integral" (Red Violet) -> "Office Theme" (Office) -
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
